# update code for mobile domain
# The backend's "mobile" error/response payloads changed shape:
#   old Unauthorized envelope: "status":401,"error":"Unauthorized"
#   new Unauthorized envelope: "code":401,"message":"Unauthorized"
#   old empty-list envelope:   []
#   new empty-list envelope:   "code":200,"message":"success","data":[]
# Update every test-data cell that encodes these fixtures across all three
# sheets, and restore each sheet's cursor/selection the way the source
# workbook has it after the edit.

$wb = $excel.ActiveWorkbook

$oldUnauthorized = '"status":401,"error":"Unauthorized"'
$newUnauthorized = '"code":401,"message":"Unauthorized"'
$newEmptyList    = '"code":200,"message":"success","data":[]'

# --- Sheet 1: "Promotion Vouchers" ---------------------------------------
$ws1 = $wb.Worksheets.Item("Promotion Vouchers")
$ws1.Range("D4").Value = $newUnauthorized
$ws1.Range("D5").Value = $newUnauthorized
$ws1.Range("D6").Value = $newEmptyList

# --- Sheet 2: "My Vouchers" ----------------------------------------------
$ws2 = $wb.Worksheets.Item("My Vouchers")
$ws2.Range("D6").Value = $newUnauthorized
$ws2.Range("D7").Value = $newUnauthorized

# --- Sheet 3: "Voucher Details" -------------------------------------------
$ws3 = $wb.Worksheets.Item("Voucher Details")
$ws3.Range("D4").Value = $newUnauthorized
$ws3.Range("D5").Value = $newUnauthorized
# These two cells also switched from General to Text number format
# (matches the "@" text style used elsewhere in column D).
$ws3.Range("D4").NumberFormat = "@"
$ws3.Range("D5").NumberFormat = "@"

# --- Restore each sheet's selection/cursor --------------------------------
$ws1.Activate()
$ws1.Range("D5").Select()

$ws2.Activate()
$ws2.Range("D7").Select()

$ws3.Activate()
$ws3.Range("D14").Select()
